# Update "想去人数" (interested-count, column F) values across the
# four sheets to reflect the latest scrape snapshot (gh-pages output
# generated at 456a3b4). Only column F values change; everything else
# (layout, other columns, formatting) stays untouched.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 17
$ws.Range("F4").Value = 1349
$ws.Range("F6").Value = 7714
$ws.Range("F8").Value = 111
$ws.Range("F10").Value = 8502
$ws.Range("F11").Value = 10
$ws.Range("F13").Value = 81
$ws.Range("F14").Value = 5699
$ws.Range("F16").Value = 2653
$ws.Range("F17").Value = 1161
$ws.Range("F18").Value = 4596
$ws.Range("F24").Value = 3654
$ws.Range("F25").Value = 73
$ws.Range("F26").Value = 42
$ws.Range("F29").Value = 3155
$ws.Range("F30").Value = 59
$ws.Range("F31").Value = 226
$ws.Range("F32").Value = 359
$ws.Range("F33").Value = 136
$ws.Range("F34").Value = 338
$ws.Range("F35").Value = 966
$ws.Range("F38").Value = 886
$ws.Range("F39").Value = 2556
$ws.Range("F43").Value = 3111
$ws.Range("F45").Value = 2300

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 109
$ws.Range("F3").Value = 135
$ws.Range("F6").Value = 6

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F3").Value = 1337

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 1337
$ws.Range("F4").Value = 17
$ws.Range("F5").Value = 1349
$ws.Range("F6").Value = 7714
$ws.Range("F7").Value = 111
$ws.Range("F9").Value = 8502
$ws.Range("F10").Value = 10
$ws.Range("F11").Value = 81
$ws.Range("F12").Value = 5699
$ws.Range("F14").Value = 2653
$ws.Range("F15").Value = 1161
$ws.Range("F16").Value = 4596
$ws.Range("F19").Value = 109
$ws.Range("F21").Value = 135
$ws.Range("F24").Value = 3654
$ws.Range("F25").Value = 73
$ws.Range("F26").Value = 42
$ws.Range("F29").Value = 3155
$ws.Range("F30").Value = 59
$ws.Range("F31").Value = 359
$ws.Range("F32").Value = 136
$ws.Range("F33").Value = 338
$ws.Range("F35").Value = 966
$ws.Range("F38").Value = 886
$ws.Range("F40").Value = 2556
$ws.Range("F44").Value = 3111
